# Apply the "BGP Interface" section insertion to the Vscs sheet.
#
# Summary of the change:
#   A new "BGP Interface" section (4 rows: a section header plus
#   "BGP Interface IP Address", "BGP Interface Prefix length" and
#   "BGP Interface VLAN ID") is inserted right before the existing
#   "vCenter Parameters" section on the "Vscs" worksheet. Everything
#   that used to live at row 33 onward moves down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

$insertAt = 33
$shift = 4

# ---------------------------------------------------------------
# 1. Remember every comment that sits on or after the insertion
#    point, because native row-insert does not relocate comments
#    together with the cells they annotate.
# ---------------------------------------------------------------
$movedComments = @()
foreach ($c in $ws.Comments) {
    $r = $c.Parent.Row
    if ($r -ge $insertAt) {
        $movedComments += , @($r, $c.Parent.Column, $c.Text())
    }
}

# Remove those comments now; they will be re-created after the shift.
foreach ($c in @($ws.Comments)) {
    $r = $c.Parent.Row
    if ($r -ge $insertAt) {
        $c.Delete()
    }
}

# ---------------------------------------------------------------
# 2. Insert the 4 new rows. This correctly shifts cell values,
#    styles, merged cells and data validations.
# ---------------------------------------------------------------
$ws.Rows("$($insertAt):$($insertAt + $shift - 1)").Insert()

# ---------------------------------------------------------------
# 3. Re-create the relocated comments at their new row.
# ---------------------------------------------------------------
foreach ($item in $movedComments) {
    $newRow = [int]$item[0] + $shift
    $col = [int]$item[1]
    $text = [string]$item[2]
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.AddComment($text)
}

# ---------------------------------------------------------------
# 4. Populate the 4 new rows with labels and formatting matching
#    the rest of the sheet.
# ---------------------------------------------------------------

# Row 33: new section header "BGP Interface" (same look as other
# section headers, e.g. the "vCenter Parameters" header which is
# now at row 37). Section header rows only carry a value/style in
# column A, so only copy that single cell's format. Merge first so
# the paste-format step does not spawn empty B/C cell records.
$ws.Range("A33:C33").Merge() | Out-Null
$ws.Range("A37").Copy()
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("A33").Value = "BGP Interface"

# Row 34: "BGP Interface IP Address" - styled like row 38 (vCenter
# Datacenter Name), a "grey" label/value row.
$ws.Range("A38:C38").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122) | Out-Null
$ws.Range("A34").Value = "BGP Interface IP Address"
$ws.Range("A34").AddComment("IP Address for Optional BGP Interface") | Out-Null

# Row 35: "BGP Interface Prefix length" - styled like row 42
# (vCenter VM Folder), a "yellow" label/value row.
$ws.Range("A42:C42").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122) | Out-Null
$ws.Range("A35").Value = "BGP Interface Prefix length"
$ws.Range("A35").AddComment("Prefix length for the optional BGP interface [default: 24]") | Out-Null

# Row 36: "BGP Interface VLAN ID" - styled like row 38 as well
# (grey label/value row).
$ws.Range("A38:C38").Copy()
$ws.Range("A36:C36").PasteSpecial(-4122) | Out-Null
$ws.Range("A36").Value = "BGP Interface VLAN ID"
$ws.Range("A36").AddComment("VLAN ID for the optional BGP interface [default: 1000]") | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 5. Data validation: prefix length and VLAN ID are integers, just
#    like the other "whole number" fields on this sheet.
# ---------------------------------------------------------------
foreach ($addr in @("B35", "C35", "B36", "C36")) {
    $rng = $ws.Range($addr)
    $rng.Validation.Add(1, 1, 3, "Please provide integer") | Out-Null
    $rng.Validation.ErrorTitle = "Invalid Entry"
    $rng.Validation.ErrorMessage = "Your entry is not an integer, change anyway?"
    $rng.Validation.PromptTitle = "Integer Selection"
    $rng.Validation.InputMessage = "Please provide integer"
    $rng.Validation.ShowInput = $true
    $rng.Validation.ShowError = $true
}

Write-Host "Done applying BGP Interface section."
